# Kitty in the Rain Course of Action - highlight status corrections
#
# "Spawn platforms in random places" (under Platforms) gets highlighted
# yellow to flag it as corrected/needing-attention, and the "Camera"
# heading gets highlighted cyan (turquoise) to match its now-complete
# "Auto Scroll same speed as cat" child item.
#
# WdColorIndex: wdYellow = 7, wdTurquoise = 3

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($text -eq "Spawn platforms in random places") {
        $p.Range.Font.HighlightColorIndex = 7
    }
    elseif ($text -eq "Camera") {
        $p.Range.Font.HighlightColorIndex = 3
    }
}

Write-Output "done"
